# Reorganize test data with sales_tax and use_tax folders
# Adds 4 new claim rows (22-25) to the synthetic Sales Tax Claims sheet,
# reflecting newly-added multi-vendor / PO-only test documents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: V0021 / Northstar Technologies Inc -----------------------
$ws.Range("A22").Value = "V0021"
$ws.Range("B22").Value = "Northstar Technologies Inc"
$ws.Range("E22").Value = "PO-Equipment-Tracker.xlsx"
$ws.Range("H22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H22").Value = [datetime]"2023-10-01"
$ws.Range("I22").Value = 154556
$ws.Range("J22").Value = 12364.48
$ws.Range("K22").Value = 166920.48
$ws.Range("L22").Value = 8
$ws.Range("M22").Value = "IT Equipment (Multi-vendor PO)"
$ws.Range("N22").Value = "Cloud Platform Access, Network Switch, Technical Consulting, Server Rack, Managed Services"
$ws.Range("P22").Value = "Equipment Tracker, Multi-vendor"

# --- Row 23: V0022 / Northstar Technologies Inc -----------------------
$ws.Range("A23").Value = "V0022"
$ws.Range("B23").Value = "Northstar Technologies Inc"
$ws.Range("E23").Value = "PO-License-Summary.xlsx"
$ws.Range("H23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H23").Value = [datetime]"2024-07-03"
$ws.Range("I23").Value = 861047
$ws.Range("J23").Value = 68883.76
$ws.Range("K23").Value = 929930.76
$ws.Range("L23").Value = 8
$ws.Range("M23").Value = "Software License (Multi-vendor PO)"
$ws.Range("N23").Value = "System Integration, Training Program, Storage Array, Professional Services, Security Suite"
$ws.Range("P23").Value = "License Summary, Multi-vendor"

# --- Row 24: V0023 / Lake Systems Group -------------------------------
$ws.Range("A24").Value = "V0023"
$ws.Range("B24").Value = "Lake Systems Group"
$ws.Range("E24").Value = "PO-Request-Email.eml"
$ws.Range("G24").Value = "PO-4900596750"
$ws.Range("H24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H24").Value = [datetime]"2024-02-23"
$ws.Range("I24").Value = 57864
$ws.Range("J24").Value = 3471.84
$ws.Range("K24").Value = 61335.84
$ws.Range("L24").Value = 6
$ws.Range("M24").Value = "IT Equipment"
$ws.Range("N24").Value = "Workstation - High Performance (5), Network Switch - 48 Port (4), UPS Battery Backup System (3)"
$ws.Range("P24").Value = "Workstation, Network Switch, UPS"

# --- Row 25: V0024 / Valley Data Systems ------------------------------
$ws.Range("A25").Value = "V0024"
$ws.Range("B25").Value = "Valley Data Systems"
$ws.Range("E25").Value = "PO-Service-Quotation.docx"
$ws.Range("G25").Value = "Q-78146"
$ws.Range("H25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H25").Value = [datetime]"2023-01-02"
$ws.Range("I25").Value = 827428
$ws.Range("J25").Value = 46335.97
$ws.Range("K25").Value = 873763.97
$ws.Range("L25").Value = 5.6
$ws.Range("M25").Value = "Professional Services"
$ws.Range("N25").Value = "Data services implementation"
$ws.Range("P25").Value = "Professional Services Quotation"
